# Adds client-side validation error codes to the "Register" sheet and
# updates sheet selections / the active tab to match the author's edit.

$wb = $excel.ActiveWorkbook

$wsAuth     = $wb.Worksheets.Item("Auth")
$wsRegister = $wb.Worksheets.Item("Register")

# --- New "Register" error-code rows -----------------------------------
# R1: username can not be empty
# R2: password can not be empty
# R3: confirm password can not be empty
# R4: passwords do not match
# R5: code can not be empty
$wsRegister.Range("B3").Value = "username can not be empty"
$wsRegister.Range("C3").Value = "ERROR"

$wsRegister.Range("B4").Value = "password can not be empty"
$wsRegister.Range("C4").Value = "ERROR"

$wsRegister.Range("B5").Value = "confirm password can not be empty"
$wsRegister.Range("C5").Value = "ERROR"

$wsRegister.Range("B6").Value = "passwords do not match"
$wsRegister.Range("C6").Value = "ERROR"

$wsRegister.Range("B7").Value = "code can not be empty"
$wsRegister.Range("C7").Value = "ERROR"

# --- Selection / active sheet updates ----------------------------------
# Auth sheet loses focus; its saved selection moves from B15 to B14.
$wsAuth.Activate()
$null = $wsAuth.Range("B14").Select()

# Register becomes the active (tabSelected) sheet with a new selection.
$wsRegister.Activate()
$null = $wsRegister.Range("G10").Select()
